# Add a new test-case row (row 26) to Sheet1, mirroring the existing
# Environment / Testdata_name / Testdata_path table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A26").Value = "test"
$ws.Range("B26").Value = "livehta_3750_manage_abbreviation"
$ws.Range("C26").Value = "\Testdata\Non_Oncology\DataFiles\ManageAbbreviations\LIVEHTA_3750_manageAbbreviation_Data.xlsx"

$ws.Range("E29").Select()
